$d = $word.ActiveDocument

# Locate the paragraph that begins "The Calculator class implementation" (it is
# split across two runs: "The Calculator class implementation " and
# "follows Python best practices and includes the following key methods:").
$anchor = $d.Content
$found = $anchor.Find.Execute("The Calculator class implementation", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
[void]$anchor.Expand(4)  # wdParagraph -> grow the range to the whole paragraph

# Replace the whole paragraph (both runs) with the new heading line. Using the
# paragraph's own start/end (rather than the Find range) makes sure both runs
# are overwritten instead of just the matched text.
$first = $d.Range($anchor.Start, $anchor.End)
$first.Text = "### 6. Technical Implementation"

# Re-fetch the (now single-run) paragraph and insert the remaining new
# paragraphs directly after it, in order.
$prevRange = $d.Range($first.Start, $first.End)

$newParagraphs = @(
    "The Calculator class implementation follows Python best practices and includes the following key methods:",
    "- **Addition (``add``)**: Performs the addition of two numbers and returns the result.",
    "- **Subtraction (``subtract``)**: Subtracts the second number from the first and returns the result.",
    "- **Multiplication (``multiply``)**: Multiplies two numbers and returns the product.",
    "- **Division (``divide``)**: Divides the first number by the second, handling division by zero with appropriate error messages.",
    "- **History Tracking**: Maintains a log of all calculations performed during the session, allowing users to review past operations.",
    "Recent updates to the codebase include enhancements to the existing methods and the addition of new functionalities:",
    "- **Enhanced Error Handling**: The error handling framework has been improved to provide more descriptive messages and handle edge cases more gracefully.",
    "- **Configuration Management**: A new configuration file (``deepdocs.yml``) has been introduced to manage application settings and environment configurations, allowing for greater flexibility and customization.",
    "- **Modular Code Structure**: The code has been refactored to improve modularity, making it easier to maintain and extend. This includes the separation of concerns into distinct modules for better organization.",
    "These updates ensure that the Calculator Application remains robust, user-friendly, and adaptable to future requirements. The implementation is designed to be easily extensible, supporting the integration of additional features and improvements over time."
)

foreach ($txt in $newParagraphs) {
    [void]$prevRange.InsertParagraphAfter()
    $prevPara = $prevRange.Paragraphs(1)
    $nextPara = $prevPara.Next()
    $nextRange = $nextPara.Range
    $nextRange.Text = $txt
    $prevRange = $d.Range($nextRange.Start, $nextRange.End)
}
